$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 12 data
$ws.Range("A12").Value = "646. Maximum Length of Pair Chain"

$ws.Range("B12").Value = "Medium"
$ws.Range("B12").Interior.Color = 49407

$ws.Range("C12").Value = "Longest Increasing Subsequence"

$ws.Range("D12").Value = "The optimal is the Greedy solution. House Robber variation. Take/Not Take. Based on the conditions, we need to skip, but if conditions are met, we need to iterate within the recursive function to find maxChain of max(maxChain, 1+solve) (take and not take for the rest of the array). The return is then the max of maxChain and solve(i+1) (not take). Remember the Java 8 lambda sort for sorting a 2d array: Arrays.sort(pairs, (a,b) -> Integer.compare(a[0], b[0])); You need memoization atleast to pass TLE."

$ws.Range("E12").Value = "https://leetcode.com/problems/maximum-length-of-pair-chain/solutions/745935/java-solution-recursion-memoziation/?envType=study-plan-v2&envId=dynamic-programming "
$ws.Hyperlinks.Add($ws.Range("E12"), "https://leetcode.com/problems/maximum-length-of-pair-chain/solutions/745935/java-solution-recursion-memoziation/?envType=study-plan-v2&envId=dynamic-programming ")
$ws.Range("E12").Style = "Hyperlink"

# Update the active selection to reflect the new state of the sheet
$ws.Range("E22").Select()
